$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")

# Add the new "Noise Algorithm" row (row 30) below the existing Canvas API row (29).
# Set the hyperlink text (column B) first, then the label (column A), so the new
# shared-string entries are created in the same order as the authored workbook.
$ws.Range("B30").Value = "https://www.npmjs.com/package/simplex-noise "
$ws.Range("A30").Value = "Noise Algorithm"

# Turn the URL in B30 into a real hyperlink, mirroring the existing B29 hyperlink.
$ws.Hyperlinks.Add($ws.Range("B30"), "https://www.npmjs.com/package/simplex-noise ") | Out-Null

# Re-apply the built-in Hyperlink style to both hyperlink cells so B30 shares the
# same cell style index as B29 instead of getting a brand new (duplicate) style.
$ws.Range("B29").Style = "Hyperlink"
$ws.Range("B30").Style = "Hyperlink"

# Widen column A to fit the new, longer label, switching it from an auto "best fit"
# width to an explicit custom width.
$ws.Columns.Item(1).ColumnWidth = 15

# Move the active cell/selection down to the newly added row, like Excel would
# leave it after typing the new entry.
$ws.Range("A30").Select() | Out-Null
